$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "PROGEN Dieseltech Services Corp."
$ws.Range("C2").Value = "Purok San Jose, Brgy. Calumangan, Bago City"
$ws.Range("H2").Value = "PAID BILLING REPORT"
$ws.Range("C3").Value = "Negros Occidental, Philippines 6101"
$ws.Range("C4").Value = "Tel. No. 476 - 7382"
